# RECO_holdings.xlsx update:
#  - Roll the "as of" date in the confidential disclaimer from 2021-05-19 to 2021-05-20.
#  - Refresh the daily Weight (col D) / Percent Change (col E) figures for rows 2-38.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected (no password needed to lift it); drop protection so the
# cells below (all locked by default) can be written, then restore it afterwards.
$ws.Unprotect()

# --- Disclaimer text (A41): bump the "as of" date by one day -----------------
$old = $ws.Range("A41").Value2
$ws.Range("A41").Value2 = $old -replace '2021-05-19', '2021-05-20'

# --- Per-row Weight (D) / Percent Change (E) values ---------------------------
# Index 0 -> row 2 ... index 35 -> row 37; $null means "leave this cell alone".
$dValues = @(
    0.0320656599695654, 0.02849729510865729, 0.02863489729953723, 0.06389407478244125,
    0.01572698143695042, 0.01549922608652845, 0.02949313165386169, 0.03462197653289532,
    0.02892611572937365, 0.02941404993496517, 0.01084605774665736, 0.01436519423755239,
    0.01422126550916073, 0.009068300705863652, 0.00816182650301232, 0.03063190840597154,
    0.02509440874454015, 0.03349308499564754, 0.03099904528594861, 0.04564418880839468,
    0.03542663302266739, 0.03108405813376238, 0.03047651282833988, 0.01494466553276662,
    0.01498203164494522, 0.0311402061541789,  0.03143340162698775, 0.02883952124718197,
    0.02891543969732263, 0.03337031062706069, 0.03162734954258145, 0.02888855191289781,
    0.03238396388812388, 0.03069715082406117, 0.0318521393285448,  0.03463937451105255,
    $null
)

$eValues = @(
    0.01087613293051359, 0.02331050846046567, 0.02400629673356947, 0.004913670400396031,
    0.0007542615779152495, 0.01515383437930495, 0.01327273458552858, 0.01815897670169031,
    0.02864465860159937, 0.02093723534393521, 0.03727670433831554, 0.02909441233140653,
    0.0173219151420787, -0.002092962413883237, -0.002398081534772167, 0.01102376434444752,
    0.003379842274027212, 0.004692757216220977, 0.01600816352562284, 0.02052228300652748,
    0.003906467994865714, 0.004611225950071463, -0.008524054180289675, 0.007977140135730476,
    0.002612826603325402, 0.004818772260632764, 0.00840293851262941, 0.01382033563672258,
    0.002331528279181594, 0.03887101647619229, 0.0123770886337069, 0.007829181494662052,
    -0.002100122100122115, 0.0004637143519590214, -0.001737943020296617, 0.01269348431578465,
    0.01107420850567364
)

$firstRow = 2
for ($i = 0; $i -lt $dValues.Count; $i++) {
    $row = $firstRow + $i
    if ($null -ne $dValues[$i]) {
        $ws.Cells.Item($row, 4).Value2 = $dValues[$i]
    }
    if ($null -ne $eValues[$i]) {
        $ws.Cells.Item($row, 5).Value2 = $eValues[$i]
    }
}

# --- Restore sheet protection (contents locked, as it was before) -------------
$ws.Protect()
